$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "wenden"
$ws.Cells.Item(2, 2).Value = "none"
$ws.Cells.Item(2, 3).Value = "none"
$ws.Cells.Item(3, 1).Value = "decken"
$ws.Cells.Item(3, 2).Value = "dog/dog027.jpg"
$ws.Cells.Item(3, 3).Value = "dog"
$ws.Cells.Item(4, 1).Value = "albern"
$ws.Cells.Item(4, 2).Value = "house/house022.jpg"
$ws.Cells.Item(4, 3).Value = "house"
$ws.Cells.Item(5, 1).Value = "rufen"
$ws.Cells.Item(5, 2).Value = "none"
$ws.Cells.Item(5, 3).Value = "none"
$ws.Cells.Item(6, 1).Value = "binden"
$ws.Cells.Item(6, 2).Value = "house/house003.jpg"
$ws.Cells.Item(6, 3).Value = "house"
$ws.Cells.Item(7, 1).Value = "drücken"
$ws.Cells.Item(7, 2).Value = "house/house031.jpg"
$ws.Cells.Item(7, 3).Value = "house"
$ws.Cells.Item(8, 1).Value = "enden"
$ws.Cells.Item(8, 2).Value = "none"
$ws.Cells.Item(8, 3).Value = "none"
$ws.Cells.Item(9, 1).Value = "klingen"
$ws.Cells.Item(9, 2).Value = "house/house014.jpg"
$ws.Cells.Item(9, 3).Value = "house"
$ws.Cells.Item(10, 1).Value = "stärken"
$ws.Cells.Item(10, 2).Value = "dog/dog023.jpg"
$ws.Cells.Item(10, 3).Value = "dog"
$ws.Cells.Item(11, 1).Value = "opfern"
$ws.Cells.Item(11, 2).Value = "none"
$ws.Cells.Item(11, 3).Value = "none"
$ws.Cells.Item(12, 1).Value = "spenden"
$ws.Cells.Item(12, 2).Value = "dog/dog006.jpg"
$ws.Cells.Item(12, 3).Value = "dog"
$ws.Cells.Item(13, 1).Value = "schwimmen"
$ws.Cells.Item(13, 2).Value = "dog/dog002.jpg"
$ws.Cells.Item(13, 3).Value = "dog"
$ws.Cells.Item(14, 1).Value = "schalten"
$ws.Cells.Item(14, 2).Value = "none"
$ws.Cells.Item(14, 3).Value = "none"
$ws.Cells.Item(15, 1).Value = "treiben"
$ws.Cells.Item(15, 2).Value = "dog/dog004.jpg"
$ws.Cells.Item(15, 3).Value = "dog"
$ws.Cells.Item(16, 1).Value = "stürmen"
$ws.Cells.Item(16, 2).Value = "dog/dog001.jpg"
$ws.Cells.Item(16, 3).Value = "dog"
$ws.Cells.Item(17, 1).Value = "klagen"
$ws.Cells.Item(17, 2).Value = "none"
$ws.Cells.Item(17, 3).Value = "none"
$ws.Cells.Item(18, 1).Value = "holen"
$ws.Cells.Item(18, 2).Value = "dog/dog031.jpg"
$ws.Cells.Item(18, 3).Value = "dog"
$ws.Cells.Item(19, 1).Value = "leugnen"
$ws.Cells.Item(19, 2).Value = "house/house019.jpg"
$ws.Cells.Item(19, 3).Value = "house"
$ws.Cells.Item(20, 1).Value = "sparen"
$ws.Cells.Item(20, 2).Value = "none"
$ws.Cells.Item(20, 3).Value = "none"
$ws.Cells.Item(21, 1).Value = "leuchten"
$ws.Cells.Item(21, 2).Value = "house/house007.jpg"
$ws.Cells.Item(21, 3).Value = "house"
$ws.Cells.Item(22, 1).Value = "hoffen"
$ws.Cells.Item(22, 2).Value = "dog/dog014.jpg"
$ws.Cells.Item(22, 3).Value = "dog"
$ws.Cells.Item(23, 1).Value = "hören"
$ws.Cells.Item(23, 2).Value = "none"
$ws.Cells.Item(23, 3).Value = "none"
$ws.Cells.Item(24, 1).Value = "deuten"
$ws.Cells.Item(24, 2).Value = "dog/dog025.jpg"
$ws.Cells.Item(24, 3).Value = "dog"
$ws.Cells.Item(25, 1).Value = "parken"
$ws.Cells.Item(25, 2).Value = "house/house008.jpg"
$ws.Cells.Item(25, 3).Value = "house"
$ws.Cells.Item(26, 1).Value = "drohen"
$ws.Cells.Item(26, 2).Value = "none"
$ws.Cells.Item(26, 3).Value = "none"
$ws.Cells.Item(27, 1).Value = "gelten"
$ws.Cells.Item(27, 2).Value = "house/house029.jpg"
$ws.Cells.Item(27, 3).Value = "house"
$ws.Cells.Item(28, 1).Value = "handeln"
$ws.Cells.Item(28, 2).Value = "house/house012.jpg"
$ws.Cells.Item(28, 3).Value = "house"
$ws.Cells.Item(29, 1).Value = "orten"
$ws.Cells.Item(29, 2).Value = "none"
$ws.Cells.Item(29, 3).Value = "none"
$ws.Cells.Item(30, 1).Value = "bitten"
$ws.Cells.Item(30, 2).Value = "dog/dog015.jpg"
$ws.Cells.Item(30, 3).Value = "dog"
$ws.Cells.Item(31, 1).Value = "regnen"
$ws.Cells.Item(31, 2).Value = "dog/dog024.jpg"
$ws.Cells.Item(31, 3).Value = "dog"
$ws.Cells.Item(32, 1).Value = "weigern"
$ws.Cells.Item(32, 2).Value = "none"
$ws.Cells.Item(32, 3).Value = "none"
$ws.Cells.Item(33, 1).Value = "starren"
$ws.Cells.Item(33, 2).Value = "dog/dog026.jpg"
$ws.Cells.Item(33, 3).Value = "dog"
$ws.Cells.Item(34, 1).Value = "ärgern"
$ws.Cells.Item(34, 2).Value = "house/house004.jpg"
$ws.Cells.Item(34, 3).Value = "house"
$ws.Cells.Item(35, 1).Value = "ächzen"
$ws.Cells.Item(35, 2).Value = "none"
$ws.Cells.Item(35, 3).Value = "none"
$ws.Cells.Item(36, 1).Value = "jubeln"
$ws.Cells.Item(36, 2).Value = "house/house018.jpg"
$ws.Cells.Item(36, 3).Value = "house"
$ws.Cells.Item(37, 1).Value = "testen"
$ws.Cells.Item(37, 2).Value = "house/house015.jpg"
$ws.Cells.Item(37, 3).Value = "house"
$ws.Cells.Item(38, 1).Value = "kosten"
$ws.Cells.Item(38, 2).Value = "none"
$ws.Cells.Item(38, 3).Value = "none"
$ws.Cells.Item(39, 1).Value = "biegen"
$ws.Cells.Item(39, 2).Value = "dog/dog013.jpg"
$ws.Cells.Item(39, 3).Value = "dog"
$ws.Cells.Item(40, 1).Value = "schreiben"
$ws.Cells.Item(40, 2).Value = "house/house013.jpg"
$ws.Cells.Item(40, 3).Value = "house"
$ws.Cells.Item(41, 1).Value = "dauern"
$ws.Cells.Item(41, 2).Value = "none"
$ws.Cells.Item(41, 3).Value = "none"
$ws.Cells.Item(42, 1).Value = "danken"
$ws.Cells.Item(42, 2).Value = "house/house027.jpg"
$ws.Cells.Item(42, 3).Value = "house"
$ws.Cells.Item(43, 1).Value = "münzen"
$ws.Cells.Item(43, 2).Value = "dog/dog021.jpg"
$ws.Cells.Item(43, 3).Value = "dog"
$ws.Cells.Item(44, 1).Value = "stören"
$ws.Cells.Item(44, 2).Value = "none"
$ws.Cells.Item(44, 3).Value = "none"
$ws.Cells.Item(45, 1).Value = "tollen"
$ws.Cells.Item(45, 2).Value = "house/house010.jpg"
$ws.Cells.Item(45, 3).Value = "house"
$ws.Cells.Item(46, 1).Value = "wagen"
$ws.Cells.Item(46, 2).Value = "dog/dog017.jpg"
$ws.Cells.Item(46, 3).Value = "dog"
$ws.Cells.Item(47, 1).Value = "bremsen"
$ws.Cells.Item(47, 2).Value = "none"
$ws.Cells.Item(47, 3).Value = "none"
$ws.Cells.Item(48, 1).Value = "retten"
$ws.Cells.Item(48, 2).Value = "house/house021.jpg"
$ws.Cells.Item(48, 3).Value = "house"
$ws.Cells.Item(49, 1).Value = "wohnen"
$ws.Cells.Item(49, 2).Value = "dog/dog007.jpg"
$ws.Cells.Item(49, 3).Value = "dog"
